$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the old "NO" column (old A), shifting everything left by one
# (old B->A "项", old C->B "结果", old D->C "参考值", old E->D "单位")
$ws.Columns("A").Delete()

# Step 2: fill in the "结果" (result) values that used to be blank placeholders,
# clean up reference-range / unit typos, and shorten long test names to single/multi
# corrected labels, per the corrector pass recorded in the commit.
# Numeric-looking results are forced to text (NumberFormat "@") so values like
# "9.0" / "1.0" keep their trailing zero instead of being normalised to 9 / 1.
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "9.0"
$ws.Range("C7").Value = "109-245"
$ws.Range("C8").Value = "24-194"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "55.6"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "37.0"
$ws.Range("C16").Value = "3.0-20.0"
$ws.Range("C18").Value = "1.7-15.0"
$ws.Range("D19").Value = "mmol/L"
$ws.Range("D21").Value = "mmol/L"
$ws.Range("D22").Value = "mmol/L"
$ws.Range("D23").Value = "mol/L"
$ws.Range("A24").Value = "载脂蛋白-B"
$ws.Range("A25").Value = "载脂蛋白-B"
$ws.Range("A26").Value = ""
$ws.Range("A27").Value = "尿素"
$ws.Range("D27").Value = "mmol/L"
$ws.Range("A28").Value = "肌酐"
$ws.Range("A29").Value = "尿酸"
$ws.Range("A30").Value = "光抑素C"
$ws.Range("A31").Value = "钾"
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "3.8"
$ws.Range("D31").Value = "mmol/L"
$ws.Range("A32").Value = "钠"
$ws.Range("D32").Value = "mmol/L"
$ws.Range("A33").Value = "氯"
$ws.Range("D33").Value = "mmol/L"
$ws.Range("A34").Value = "总二氧化碳"
$ws.Range("C34").Value = "20.0-29.0"
$ws.Range("D34").Value = "mmol/L"
$ws.Range("A35").Value = "钙"
$ws.Range("D35").Value = "mmol/L"
$ws.Range("A36").Value = "磷"
$ws.Range("D36").Value = "mmol/L"
$ws.Range("A37").Value = ""
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "7.2"
$ws.Range("D37").Value = "g/L"
$ws.Range("A38").Value = "球蛋白"
$ws.Range("C38").Value = "0.72-4.29"
$ws.Range("D38").Value = "g/L"
$ws.Range("A39").Value = "球蛋白"
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "1.0"
$ws.Range("C39").Value = "0.29-3.44"
$ws.Range("D39").Value = "g/L"
$ws.Range("A40").Value = "补体C3"
$ws.Range("D40").Value = "mg/dL"
$ws.Range("A41").Value = "补体C4"
$ws.Range("A42").Value = "睡液酸"
$ws.Range("A43").Value = "阴离子隙"
$ws.Range("C43").Value = "8.0-16.0"
$ws.Range("D43").Value = "mmol/L"
$ws.Range("A44").Value = "同型半胱氨酸"
